$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update final score (K) values
$ws.Range("K2").Value = 59.5
$ws.Range("K3").Value = 55.5
$ws.Range("K4").Value = 53.5
$ws.Range("K5").Value = 52.5

# Update MACRO_SCORE (N) values
$ws.Range("N2").Value = 51.53902399942638
$ws.Range("N3").Value = 51.53902399942638
$ws.Range("N4").Value = 51.53902399942638
$ws.Range("N5").Value = 51.53902399942638

# Update 판단 (M2) from buy-watch message to hold/wait message
$ws.Range("M2").Value = "⛔ 관망하십시오."
